# The workbook gains a header label "Category" in A1 (using the same bold /
# bordered / centered header style already applied to B1:W1), and the
# per-cell header style that had been mistakenly left on A2:A46 (the
# category-name column) is cleared back to the default/no style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bring A1 up to the same formatting as the rest of row 1's header cells,
# then give it its text.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("A1").Value = "Category"

# A2:A46 previously carried that same header style (s="1"); drop it so the
# cells fall back to the default style, matching the rest of the data rows.
$ws.Range("A2:A46").ClearFormats()
